$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ListSheet")

# Row 2: fill in actual student data (name/lastname) and bump degree numeral (K2) 1 -> 2
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "Lastname"
$ws.Range("K2").Value = 2

# Row 3: fill in actual student data (name split across B3/C3), update gender (I3) and degree letter (J3)
$ws.Range("B3").Value = "Fem"
$ws.Range("C3").Value = "ale"
$ws.Range("I3").Value = "female"
$ws.Range("J3").Value = "A"

# Remove the stray leftover row 11 (only contained an empty formatted cell E11)
$ws.Rows.Item(11).Delete()

# Update the active selection to reflect the last edited cell
$ws.Activate()
$ws.Range("J3").Select() | Out-Null
